$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# --- Update year-range header labels (shift one year forward: drop 1396/12, add 1401/12) ---
Set-Cell 8 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 8 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 8 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 8 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 8 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 27 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 27 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 27 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 27 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 27 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 34 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 34 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 34 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 34 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 34 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 41 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 41 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 41 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 41 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 41 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 48 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 48 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 48 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 48 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 48 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 55 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 55 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 55 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 55 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 55 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 62 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 62 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 62 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 62 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 62 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 69 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 69 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 69 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 69 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 69 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 76 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 76 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 76 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 76 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 76 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 83 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 83 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 83 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 83 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 83 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 89 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 89 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 89 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 89 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 89 9 "دوازده ماهه منتهی به 1401/12"
Set-Cell 95 5 "دوازده ماهه منتهی به 1397/12"
Set-Cell 95 6 "دوازده ماهه منتهی به 1398/12"
Set-Cell 95 7 "دوازده ماهه منتهی به 1399/12"
Set-Cell 95 8 "دوازده ماهه منتهی به 1400/12"
Set-Cell 95 9 "دوازده ماهه منتهی به 1401/12"

# --- Update yearly data columns (shift E<-F<-G<-H<-I, with new figures for the latest year 1401/12) ---
Set-Cell 10 5 2214862
Set-Cell 10 6 2734511
Set-Cell 10 7 4739492
Set-Cell 10 8 15018423
Set-Cell 10 9 13279766
Set-Cell 11 5 122826
Set-Cell 11 6 163928
Set-Cell 11 7 235590
Set-Cell 11 8 341234
Set-Cell 11 9 554944
Set-Cell 12 5 1765091
Set-Cell 12 6 2494937
Set-Cell 12 7 3197294
Set-Cell 12 8 7241278
Set-Cell 12 9 12087404
Set-Cell 13 5 4102779
Set-Cell 13 6 5393376
Set-Cell 13 7 8172376
Set-Cell 13 8 22600935
Set-Cell 13 9 25922114
Set-Cell 14 5 0
Set-Cell 14 6 -41480
Set-Cell 14 7 -3002
Set-Cell 14 8 -11324
Set-Cell 14 9 -60782
Set-Cell 15 5 4102779
Set-Cell 15 6 5351896
Set-Cell 15 7 8169374
Set-Cell 15 8 22589611
Set-Cell 15 9 25861332
Set-Cell 16 5 0
Set-Cell 16 6 0
Set-Cell 16 7 0
Set-Cell 16 8 0
Set-Cell 16 9 0
Set-Cell 17 5 0
Set-Cell 17 6 0
Set-Cell 17 7 0
Set-Cell 17 8 0
Set-Cell 17 9 -21058
Set-Cell 18 5 4102779
Set-Cell 18 6 5351896
Set-Cell 18 7 8169374
Set-Cell 18 8 22589611
Set-Cell 18 9 25840274
Set-Cell 19 5 250893
Set-Cell 19 6 516698
Set-Cell 19 7 150145
Set-Cell 19 8 715957
Set-Cell 19 9 1850267
Set-Cell 20 5 -525702
Set-Cell 20 6 -150146
Set-Cell 20 7 -715957
Set-Cell 20 8 -1850267
Set-Cell 20 9 -2217271
Set-Cell 21 5 3827970
Set-Cell 21 6 5718448
Set-Cell 21 7 7603562
Set-Cell 21 8 21455301
Set-Cell 21 9 25473270
Set-Cell 22 5 0
Set-Cell 22 6 0
Set-Cell 22 7 0
Set-Cell 22 8 0
Set-Cell 22 9 0
Set-Cell 23 5 3827970
Set-Cell 23 6 5718448
Set-Cell 23 7 7603562
Set-Cell 23 8 21455301
Set-Cell 23 9 25473270
Set-Cell 29 5 0
Set-Cell 29 6 0
Set-Cell 29 7 0
Set-Cell 29 8 0
Set-Cell 29 9 0
Set-Cell 30 5 0
Set-Cell 30 6 0
Set-Cell 30 7 0
Set-Cell 30 8 0
Set-Cell 30 9 0
Set-Cell 36 5 294929189
Set-Cell 36 6 271140880
Set-Cell 36 7 293215587
Set-Cell 36 8 288892075
Set-Cell 36 9 251592030
Set-Cell 37 5 294929189
Set-Cell 37 6 271140880
Set-Cell 37 7 293215587
Set-Cell 37 8 288892075
Set-Cell 37 9 251592030
Set-Cell 43 5 294929189
Set-Cell 43 6 271140880
Set-Cell 43 7 293215587
Set-Cell 43 8 288892075
Set-Cell 43 9 251592030
Set-Cell 44 5 294929189
Set-Cell 44 6 271140880
Set-Cell 44 7 293215587
Set-Cell 44 8 288892075
Set-Cell 44 9 251592030
Set-Cell 50 5 0
Set-Cell 50 6 0
Set-Cell 50 7 0
Set-Cell 50 8 0
Set-Cell 50 9 0
Set-Cell 51 5 0
Set-Cell 51 6 0
Set-Cell 51 7 0
Set-Cell 51 8 0
Set-Cell 51 9 0
Set-Cell 57 5 0
Set-Cell 57 6 0
Set-Cell 57 7 0
Set-Cell 57 8 0
Set-Cell 57 9 0
Set-Cell 58 5 0
Set-Cell 58 6 0
Set-Cell 58 7 0
Set-Cell 58 8 0
Set-Cell 58 9 0
Set-Cell 64 5 2214862
Set-Cell 64 6 2734511
Set-Cell 64 7 4739492
Set-Cell 64 8 15018423
Set-Cell 64 9 13279766
Set-Cell 65 5 2214862
Set-Cell 65 6 2734511
Set-Cell 65 7 4739492
Set-Cell 65 8 15018423
Set-Cell 65 9 13279766
Set-Cell 71 5 2214862
Set-Cell 71 6 2734511
Set-Cell 71 7 4739492
Set-Cell 71 8 15018423
Set-Cell 71 9 13279766
Set-Cell 72 5 2214862
Set-Cell 72 6 2734511
Set-Cell 72 7 4739492
Set-Cell 72 8 15018423
Set-Cell 72 9 13279766
Set-Cell 78 5 0
Set-Cell 78 6 0
Set-Cell 78 7 0
Set-Cell 78 8 0
Set-Cell 78 9 0
Set-Cell 79 5 0
Set-Cell 79 6 0
Set-Cell 79 7 0
Set-Cell 79 8 0
Set-Cell 79 9 0
Set-Cell 85 5 7510
Set-Cell 85 6 10085
Set-Cell 85 7 16164
Set-Cell 85 8 51986
Set-Cell 85 9 52783
Set-Cell 91 5 7510
Set-Cell 91 6 10085
Set-Cell 91 7 16164
Set-Cell 91 8 51986
Set-Cell 91 9 52783
Set-Cell 97 5 0
Set-Cell 97 6 0
Set-Cell 97 7 0
Set-Cell 97 8 0
Set-Cell 97 9 0
Set-Cell 98 5 0
Set-Cell 98 6 0
Set-Cell 98 7 0
Set-Cell 98 8 0
Set-Cell 98 9 0
Set-Cell 99 5 0
Set-Cell 99 6 0
Set-Cell 99 7 0
Set-Cell 99 8 0
Set-Cell 99 9 0
Set-Cell 100 5 0
Set-Cell 100 6 0
Set-Cell 100 7 0
Set-Cell 100 8 0
Set-Cell 100 9 0
Set-Cell 101 5 176251
Set-Cell 101 6 267535
Set-Cell 101 7 329297
Set-Cell 101 8 688758
Set-Cell 101 9 869416
Set-Cell 102 5 338661
Set-Cell 102 6 607834
Set-Cell 102 7 663558
Set-Cell 102 8 3439252
Set-Cell 102 9 6266564
Set-Cell 103 5 360907
Set-Cell 103 6 381428
Set-Cell 103 7 399618
Set-Cell 103 8 443342
Set-Cell 103 9 515377
Set-Cell 104 5 449404
Set-Cell 104 6 607660
Set-Cell 104 7 876474
Set-Cell 104 8 1378123
Set-Cell 104 9 2224096
Set-Cell 105 5 0
Set-Cell 105 6 0
Set-Cell 105 7 0
Set-Cell 105 8 0
Set-Cell 105 9 0
Set-Cell 106 5 439868
Set-Cell 106 6 630480
Set-Cell 106 7 928347
Set-Cell 106 8 1291803
Set-Cell 106 9 2211951
Set-Cell 107 5 1765091
Set-Cell 107 6 2494937
Set-Cell 107 7 3197294
Set-Cell 107 8 7241278
Set-Cell 107 9 12087404

